$wb = $excel.ActiveWorkbook

# Rename first sheet: "ExcelModuleDemoToDoItem" -> "DemoToDoItem"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "DemoToDoItem"

$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Re-fit row heights on every sheet so previously hard-coded ht="15" rows
# fall back to the (now-current) default row height instead of carrying an
# explicit/custom row height.
$ws1.Range("A1:E13").EntireRow.AutoFit()
$ws2.Range("A1:C8").EntireRow.AutoFit()
$ws3.Range("A1:C6").EntireRow.AutoFit()

# Make the first sheet the active / selected tab (previously Sheet3 was
# active/selected).
$ws1.Activate()
